# "Generate Report for Handback"
# The handback CI run finished processing the last (cbf32fc4-...) file for
# both locales, so its "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) on the final row (row 14) of each
# locale sheet get stamped with the new handoff/handback timestamps.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D14").Value = "2016-03-08 06:24:29"
$zhcn.Range("G14").Value = "2016-03-08 06:24:47"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D14").Value = "2016-03-08 06:24:32"
$dede.Range("G14").Value = "2016-03-08 06:24:52"
